# Swap the per-match data (everything except the leading id column A)
# between pairs of rows. The row "id" in column A stays put while all
# other fields (match id, teams, odds, etc.) trade places between the
# two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

$rowPairs = @(
    @(105, 106),
    @(107, 108),
    @(132, 133)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valueA = $cellA.Value()
        $valueB = $cellB.Value()

        $cellA.Value = $valueB
        $cellB.Value = $valueA
    }
}
